# Project - 2 iteration
# Uppercase the header row, and add operation-routing data rows 2-11
# (columns K/L/M "Тпз"/"Тшт"/"КОИД" plus column H "№ операции" for the
# new OP80/OP90 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers switched to upper case -------------------------------
$ws.Range("A1").Value = "ОБОЗНАЧЕНИЕ"
$ws.Range("B1").Value = "НАИМЕНОВАНИЕ"
$ws.Range("C1").Value = "МАРШРУТ"
$ws.Range("D1").Value = "ВХОДИМОСТЬ"
$ws.Range("E1").Value = "ПАРТИЯ"
$ws.Range("F1").Value = "ЦЕНА за шт."
$ws.Range("G1").Value = "ЦЕНА за комплект"
$ws.Range("H1").Value = "№ ОПЕРАЦИИ"
$ws.Range("I1").Value = "НАИМЕНОВАНИЕ ОПЕРАЦИИ"
$ws.Range("J1").Value = "ОБОРУДОВАНИЕ"

# --- Helper: write a value as TEXT, even if it looks like a number -------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Row 2: append K/L/M to the existing data row -------------------------
Set-TextValue $ws.Range("K2") "10"
$ws.Range("L2").Value = "5,5"
$ws.Range("M2").Value = "QWP"

# --- Row 3 ------------------------------------------------------------------
Set-TextValue $ws.Range("K3") "10"
$ws.Range("L3").Value = "1,0"
$ws.Range("M3").Value = "QWP"

# --- Row 4 ------------------------------------------------------------------
Set-TextValue $ws.Range("K4") "120"
Set-TextValue $ws.Range("L4") "18"
Set-TextValue $ws.Range("M4") "6"

# --- Row 5 ------------------------------------------------------------------
Set-TextValue $ws.Range("K5") "100"
Set-TextValue $ws.Range("L5") "42"
Set-TextValue $ws.Range("M5") "6"

# --- Row 6 ------------------------------------------------------------------
Set-TextValue $ws.Range("K6") "0"
$ws.Range("L6").Value = "6,0"
Set-TextValue $ws.Range("M6") "50"

# --- Row 7 ------------------------------------------------------------------
Set-TextValue $ws.Range("K7") "0"
$ws.Range("L7").Value = "0,25"
Set-TextValue $ws.Range("M7") "1"

# --- Row 8 ------------------------------------------------------------------
$ws.Range("H8").Value = "OP80"
Set-TextValue $ws.Range("K8") "10"
Set-TextValue $ws.Range("L8") "10"

# --- Row 9 ------------------------------------------------------------------
$ws.Range("H9").Value = "OP90"
Set-TextValue $ws.Range("K9") "10"
Set-TextValue $ws.Range("L9") "10"

# --- Row 10 -----------------------------------------------------------------
$ws.Range("H10").Value = "OP80"
Set-TextValue $ws.Range("K10") "10"
Set-TextValue $ws.Range("L10") "10"

# --- Row 11 -----------------------------------------------------------------
$ws.Range("H11").Value = "OP90"
Set-TextValue $ws.Range("K11") "10"
Set-TextValue $ws.Range("L11") "10"
